$wb = $excel.ActiveWorkbook

$edits = @(
    @{Sheet="ALC"; Row=2; Col=8; Value=4196.9},
    @{Sheet="ALC"; Row=2; Col=9; Value=798.8},
    @{Sheet="ALC"; Row=2; Col=11; Value=798.8},
    @{Sheet="ALC"; Row=2; Col=13; Value=-685.8},
    @{Sheet="ALC"; Row=4; Col=8; Value=1656.125},
    @{Sheet="ALC"; Row=4; Col=9; Value=1656.125},
    @{Sheet="ALC"; Row=4; Col=11; Value=1656.125},
    @{Sheet="ALC"; Row=4; Col=13; Value=-1542.125},
    @{Sheet="ALC"; Row=18; Col=8; Value=1975},
    @{Sheet="ALC"; Row=18; Col=9; Value=1966.6666},
    @{Sheet="ALC"; Row=18; Col=11; Value=1966.6666},
    @{Sheet="ALC"; Row=18; Col=13; Value=-1682.6666},
    @{Sheet="ALC"; Row=29; Col=8; Value=301.16666},
    @{Sheet="ALC"; Row=29; Col=10; Value=359.4},
    @{Sheet="ALC"; Row=29; Col=12; Value=1078.2},
    @{Sheet="ALC"; Row=29; Col=14; Value=-1640.2},
    @{Sheet="ALC"; Row=40; Col=8; Value=2166.5},
    @{Sheet="ALC"; Row=40; Col=9; Value=2136.889},
    @{Sheet="ALC"; Row=40; Col=10; Value=2219.8},
    @{Sheet="ALC"; Row=40; Col=11; Value=2136.889},
    @{Sheet="ALC"; Row=40; Col=12; Value=2219.8},
    @{Sheet="ALC"; Row=40; Col=13; Value=-1961.889},
    @{Sheet="ALC"; Row=40; Col=14; Value=-2569.8},
    @{Sheet="ALC"; Row=47; Col=8; Value=13755.333},
    @{Sheet="ALC"; Row=47; Col=9; Value=8133.5},
    @{Sheet="ALC"; Row=47; Col=10; Value=24999},
    @{Sheet="ALC"; Row=47; Col=11; Value=8133.5},
    @{Sheet="ALC"; Row=47; Col=12; Value=24999},
    @{Sheet="ALC"; Row=47; Col=13; Value=-7161.5},
    @{Sheet="ALC"; Row=47; Col=14; Value=-26943},
    @{Sheet="ALC"; Row=55; Col=8; Value=762},
    @{Sheet="ALC"; Row=55; Col=9; Value=480},
    @{Sheet="ALC"; Row=55; Col=11; Value=480},
    @{Sheet="ALC"; Row=55; Col=13; Value=-266},
    @{Sheet="ALC"; Row=98; Col=8; Value=1118.3846},
    @{Sheet="ALC"; Row=98; Col=9; Value=734},
    @{Sheet="ALC"; Row=98; Col=10; Value=2399.6667},
    @{Sheet="ALC"; Row=98; Col=11; Value=734},
    @{Sheet="ALC"; Row=98; Col=12; Value=2399.6667},
    @{Sheet="ALC"; Row=98; Col=13; Value=764},
    @{Sheet="ALC"; Row=98; Col=14; Value=-5395.6667},
    @{Sheet="ALC"; Row=122; Col=8; Value=1118.3846},
    @{Sheet="ALC"; Row=122; Col=9; Value=734},
    @{Sheet="ALC"; Row=122; Col=10; Value=2399.6667},
    @{Sheet="ALC"; Row=122; Col=11; Value=2202},
    @{Sheet="ALC"; Row=122; Col=12; Value=7199.000100000001},
    @{Sheet="ALC"; Row=122; Col=13; Value=248},
    @{Sheet="ALC"; Row=122; Col=14; Value=-12099.0001},
    @{Sheet="ALC"; Row=132; Col=8; Value=1393.5},
    @{Sheet="ALC"; Row=132; Col=9; Value=1250.7858},
    @{Sheet="ALC"; Row=132; Col=11; Value=3752.3574},
    @{Sheet="ALC"; Row=132; Col=13; Value=-1222.3574},
    @{Sheet="ALC"; Row=137; Col=8; Value=2376.7273},
    @{Sheet="ALC"; Row=137; Col=10; Value=2248.8},
    @{Sheet="ALC"; Row=137; Col=12; Value=6746.400000000001},
    @{Sheet="ALC"; Row=137; Col=14; Value=-11846.4},
    @{Sheet="ALC"; Row=138; Col=8; Value=3394.38},
    @{Sheet="ALC"; Row=138; Col=9; Value=3820.125},
    @{Sheet="ALC"; Row=138; Col=10; Value=3194.0293},
    @{Sheet="ALC"; Row=138; Col=11; Value=11460.375},
    @{Sheet="ALC"; Row=138; Col=12; Value=9582.0879},
    @{Sheet="ALC"; Row=138; Col=13; Value=-6320.375},
    @{Sheet="ALC"; Row=138; Col=14; Value=-19862.0879},
    @{Sheet="ARM"; Row=2; Col=8; Value=4443.1},
    @{Sheet="ARM"; Row=2; Col=9; Value=5107.6},
    @{Sheet="ARM"; Row=2; Col=11; Value=5107.6},
    @{Sheet="ARM"; Row=2; Col=13; Value=-4994.6},
    @{Sheet="ARM"; Row=5; Col=8; Value=276.86667},
    @{Sheet="ARM"; Row=5; Col=10; Value=154},
    @{Sheet="ARM"; Row=5; Col=12; Value=154},
    @{Sheet="ARM"; Row=5; Col=14; Value=-378},
    @{Sheet="ARM"; Row=63; Col=8; Value=2999.5},
    @{Sheet="ARM"; Row=63; Col=9; Value=2999.5},
    @{Sheet="ARM"; Row=63; Col=11; Value=2999.5},
    @{Sheet="ARM"; Row=63; Col=13; Value=-2313.5},
    @{Sheet="ARM"; Row=66; Col=8; Value=2999.5},
    @{Sheet="ARM"; Row=66; Col=9; Value=2999.5},
    @{Sheet="ARM"; Row=66; Col=11; Value=14997.5},
    @{Sheet="ARM"; Row=66; Col=13; Value=-11565.5},
    @{Sheet="ARM"; Row=102; Col=8; Value=1672.1333},
    @{Sheet="ARM"; Row=102; Col=9; Value=1800.0769},
    @{Sheet="ARM"; Row=102; Col=11; Value=1800.0769},
    @{Sheet="ARM"; Row=102; Col=13; Value=-178.0769},
    @{Sheet="ARM"; Row=116; Col=8; Value=4443.1},
    @{Sheet="ARM"; Row=116; Col=9; Value=5107.6},
    @{Sheet="ARM"; Row=116; Col=11; Value=5107.6},
    @{Sheet="ARM"; Row=116; Col=13; Value=-2813.6},
    @{Sheet="ARM"; Row=132; Col=8; Value=2894.4614},
    @{Sheet="ARM"; Row=132; Col=9; Value=2894.4614},
    @{Sheet="ARM"; Row=132; Col=10; Value=0},
    @{Sheet="ARM"; Row=132; Col=11; Value=8683.3842},
    @{Sheet="ARM"; Row=132; Col=12; Value=0},
    @{Sheet="ARM"; Row=132; Col=13; Value=-6153.3842},
    @{Sheet="ARM"; Row=132; Col=14; Value=$null},
    @{Sheet="ARM"; Row=135; Col=8; Value=228329.5},
    @{Sheet="ARM"; Row=135; Col=10; Value=228329.5},
    @{Sheet="ARM"; Row=135; Col=12; Value=228329.5},
    @{Sheet="ARM"; Row=135; Col=14; Value=-238469.5},
    @{Sheet="BSM"; Row=3; Col=8; Value=4443.1},
    @{Sheet="BSM"; Row=3; Col=9; Value=5107.6},
    @{Sheet="BSM"; Row=3; Col=11; Value=5107.6},
    @{Sheet="BSM"; Row=3; Col=13; Value=-4993.6},
    @{Sheet="BSM"; Row=4; Col=8; Value=276.86667},
    @{Sheet="BSM"; Row=4; Col=10; Value=154},
    @{Sheet="BSM"; Row=4; Col=12; Value=154},
    @{Sheet="BSM"; Row=4; Col=14; Value=-384},
    @{Sheet="BSM"; Row=22; Col=8; Value=1125.0834},
    @{Sheet="BSM"; Row=22; Col=9; Value=1088.25},
    @{Sheet="BSM"; Row=22; Col=10; Value=1198.75},
    @{Sheet="BSM"; Row=22; Col=11; Value=1088.25},
    @{Sheet="BSM"; Row=22; Col=12; Value=1198.75},
    @{Sheet="BSM"; Row=22; Col=13; Value=-915.25},
    @{Sheet="BSM"; Row=22; Col=14; Value=-1544.75},
    @{Sheet="BSM"; Row=105; Col=8; Value=2270.3333},
    @{Sheet="BSM"; Row=105; Col=9; Value=1724.4},
    @{Sheet="BSM"; Row=105; Col=11; Value=1724.4},
    @{Sheet="BSM"; Row=105; Col=13; Value=22.59999999999991},
    @{Sheet="CRP"; Row=7; Col=8; Value=51},
    @{Sheet="CRP"; Row=7; Col=9; Value=51.0625},
    @{Sheet="CRP"; Row=7; Col=10; Value=50.666668},
    @{Sheet="CRP"; Row=7; Col=11; Value=51.0625},
    @{Sheet="CRP"; Row=7; Col=12; Value=50.666668},
    @{Sheet="CRP"; Row=7; Col=13; Value=61.9375},
    @{Sheet="CRP"; Row=7; Col=14; Value=-276.666668},
    @{Sheet="CRP"; Row=58; Col=8; Value=1366.3334},
    @{Sheet="CRP"; Row=58; Col=9; Value=1349.5},
    @{Sheet="CRP"; Row=58; Col=10; Value=1400},
    @{Sheet="CRP"; Row=58; Col=11; Value=1349.5},
    @{Sheet="CRP"; Row=58; Col=12; Value=1400},
    @{Sheet="CRP"; Row=58; Col=13; Value=-1146.5},
    @{Sheet="CRP"; Row=58; Col=14; Value=-1806},
    @{Sheet="CRP"; Row=105; Col=8; Value=3789.8},
    @{Sheet="CRP"; Row=105; Col=9; Value=3209.6},
    @{Sheet="CRP"; Row=105; Col=11; Value=3209.6},
    @{Sheet="CRP"; Row=105; Col=13; Value=-1462.6},
    @{Sheet="CRP"; Row=107; Col=8; Value=876.8570999999999},
    @{Sheet="CRP"; Row=107; Col=9; Value=444.55554},
    @{Sheet="CRP"; Row=107; Col=11; Value=444.55554},
    @{Sheet="CRP"; Row=107; Col=13; Value=1475.44446},
    @{Sheet="CRP"; Row=132; Col=8; Value=1979},
    @{Sheet="CRP"; Row=132; Col=9; Value=1374.5},
    @{Sheet="CRP"; Row=132; Col=10; Value=2785},
    @{Sheet="CRP"; Row=132; Col=11; Value=4123.5},
    @{Sheet="CRP"; Row=132; Col=12; Value=8355},
    @{Sheet="CRP"; Row=132; Col=13; Value=-1593.5},
    @{Sheet="CRP"; Row=132; Col=14; Value=-13415},
    @{Sheet="CRP"; Row=136; Col=8; Value=1366.3334},
    @{Sheet="CRP"; Row=136; Col=9; Value=1349.5},
    @{Sheet="CRP"; Row=136; Col=10; Value=1400},
    @{Sheet="CRP"; Row=136; Col=11; Value=4048.5},
    @{Sheet="CRP"; Row=136; Col=12; Value=4200},
    @{Sheet="CRP"; Row=136; Col=13; Value=-1498.5},
    @{Sheet="CRP"; Row=136; Col=14; Value=-9300},
    @{Sheet="CRP"; Row=141; Col=8; Value=254144.69},
    @{Sheet="CRP"; Row=141; Col=10; Value=254144.69},
    @{Sheet="CRP"; Row=141; Col=12; Value=254144.69},
    @{Sheet="CRP"; Row=141; Col=14; Value=-264504.69},
    @{Sheet="CUL"; Row=8; Col=8; Value=987.6667},
    @{Sheet="CUL"; Row=8; Col=9; Value=987.6667},
    @{Sheet="CUL"; Row=8; Col=11; Value=2963.0001},
    @{Sheet="CUL"; Row=8; Col=13; Value=-2824.0001},
    @{Sheet="CUL"; Row=97; Col=8; Value=5683154.5},
    @{Sheet="CUL"; Row=97; Col=10; Value=6945855.5},
    @{Sheet="CUL"; Row=97; Col=12; Value=20837566.5},
    @{Sheet="CUL"; Row=97; Col=14; Value=-20838558.5},
    @{Sheet="CUL"; Row=132; Col=8; Value=5713.857},
    @{Sheet="CUL"; Row=132; Col=10; Value=5000},
    @{Sheet="CUL"; Row=132; Col=12; Value=45000},
    @{Sheet="CUL"; Row=132; Col=14; Value=-50060},
    @{Sheet="GSM"; Row=43; Col=8; Value=1500},
    @{Sheet="GSM"; Row=43; Col=9; Value=1500},
    @{Sheet="GSM"; Row=43; Col=11; Value=1500},
    @{Sheet="GSM"; Row=43; Col=13; Value=-1349},
    @{Sheet="GSM"; Row=46; Col=8; Value=94833.336},
    @{Sheet="GSM"; Row=46; Col=10; Value=0},
    @{Sheet="GSM"; Row=46; Col=12; Value=0},
    @{Sheet="GSM"; Row=46; Col=14; Value=$null},
    @{Sheet="GSM"; Row=113; Col=8; Value=5005.5},
    @{Sheet="GSM"; Row=113; Col=9; Value=5005.5},
    @{Sheet="GSM"; Row=113; Col=11; Value=5005.5},
    @{Sheet="GSM"; Row=113; Col=13; Value=-2835.5},
    @{Sheet="GSM"; Row=122; Col=8; Value=3305},
    @{Sheet="GSM"; Row=122; Col=9; Value=3446.4},
    @{Sheet="GSM"; Row=122; Col=11; Value=10339.2},
    @{Sheet="GSM"; Row=122; Col=13; Value=-7889.200000000001},
    @{Sheet="GSM"; Row=132; Col=8; Value=1918.8},
    @{Sheet="GSM"; Row=132; Col=9; Value=1918.8},
    @{Sheet="GSM"; Row=132; Col=10; Value=0},
    @{Sheet="GSM"; Row=132; Col=11; Value=5756.4},
    @{Sheet="GSM"; Row=132; Col=12; Value=0},
    @{Sheet="GSM"; Row=132; Col=13; Value=-3226.4},
    @{Sheet="GSM"; Row=132; Col=14; Value=$null},
    @{Sheet="LTW"; Row=16; Col=8; Value=86.666664},
    @{Sheet="LTW"; Row=16; Col=9; Value=86.666664},
    @{Sheet="LTW"; Row=16; Col=11; Value=86.666664},
    @{Sheet="LTW"; Row=16; Col=13; Value=83.333336},
    @{Sheet="LTW"; Row=43; Col=8; Value=0},
    @{Sheet="LTW"; Row=43; Col=9; Value=0},
    @{Sheet="LTW"; Row=43; Col=10; Value=0},
    @{Sheet="LTW"; Row=43; Col=11; Value=0},
    @{Sheet="LTW"; Row=43; Col=12; Value=0},
    @{Sheet="LTW"; Row=43; Col=13; Value=$null},
    @{Sheet="LTW"; Row=43; Col=14; Value=$null},
    @{Sheet="LTW"; Row=132; Col=8; Value=3212.8865},
    @{Sheet="LTW"; Row=132; Col=9; Value=3011.439},
    @{Sheet="LTW"; Row=132; Col=10; Value=5966},
    @{Sheet="LTW"; Row=132; Col=11; Value=9034.316999999999},
    @{Sheet="LTW"; Row=132; Col=12; Value=17898},
    @{Sheet="LTW"; Row=132; Col=13; Value=-6504.316999999999},
    @{Sheet="LTW"; Row=132; Col=14; Value=-22958},
    @{Sheet="WVR"; Row=10; Col=8; Value=1005},
    @{Sheet="WVR"; Row=10; Col=9; Value=1005},
    @{Sheet="WVR"; Row=10; Col=10; Value=0},
    @{Sheet="WVR"; Row=10; Col=11; Value=1005},
    @{Sheet="WVR"; Row=10; Col=12; Value=0},
    @{Sheet="WVR"; Row=10; Col=13; Value=-836},
    @{Sheet="WVR"; Row=10; Col=14; Value=$null},
    @{Sheet="WVR"; Row=14; Col=8; Value=1651},
    @{Sheet="WVR"; Row=14; Col=9; Value=1330},
    @{Sheet="WVR"; Row=14; Col=10; Value=1865},
    @{Sheet="WVR"; Row=14; Col=11; Value=1330},
    @{Sheet="WVR"; Row=14; Col=12; Value=1865},
    @{Sheet="WVR"; Row=14; Col=13; Value=-1162},
    @{Sheet="WVR"; Row=14; Col=14; Value=-2201},
    @{Sheet="WVR"; Row=70; Col=8; Value=54547.5},
    @{Sheet="WVR"; Row=70; Col=9; Value=54547.5},
    @{Sheet="WVR"; Row=70; Col=11; Value=54547.5},
    @{Sheet="WVR"; Row=70; Col=13; Value=-54232.5},
    @{Sheet="WVR"; Row=73; Col=8; Value=54547.5},
    @{Sheet="WVR"; Row=73; Col=9; Value=54547.5},
    @{Sheet="WVR"; Row=73; Col=11; Value=54547.5},
    @{Sheet="WVR"; Row=73; Col=13; Value=-53455.5},
    @{Sheet="WVR"; Row=126; Col=8; Value=4417.875},
    @{Sheet="WVR"; Row=126; Col=9; Value=3064.875},
    @{Sheet="WVR"; Row=126; Col=10; Value=7123.875},
    @{Sheet="WVR"; Row=126; Col=11; Value=9194.625},
    @{Sheet="WVR"; Row=126; Col=12; Value=21371.625},
    @{Sheet="WVR"; Row=126; Col=13; Value=-6724.625},
    @{Sheet="WVR"; Row=126; Col=14; Value=-26311.625},
    @{Sheet="WVR"; Row=136; Col=8; Value=2415.2307},
    @{Sheet="WVR"; Row=136; Col=9; Value=2616.6},
    @{Sheet="WVR"; Row=136; Col=10; Value=1744},
    @{Sheet="WVR"; Row=136; Col=11; Value=7849.799999999999},
    @{Sheet="WVR"; Row=136; Col=12; Value=5232},
    @{Sheet="WVR"; Row=136; Col=13; Value=-5299.799999999999},
    @{Sheet="WVR"; Row=136; Col=14; Value=-10332}
)

foreach ($edit in $edits) {
    $ws = $wb.Worksheets.Item($edit.Sheet)
    $cell = $ws.Cells.Item($edit.Row, $edit.Col)
    if ($null -eq $edit.Value) {
        $cell.ClearContents()
    } else {
        $cell.Value = $edit.Value
    }
}
